$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Metaphor/Simile"
$ws.Range("F1").Value = "Stories / anecdotes"
$ws.Range("H1").Value = "Lists / Repetition "
$ws.Range("I1").Value = "Moral conviction"
$ws.Range("J1").Value = "Sentiment of the collective"
$ws.Range("K1").Value = "Ambitious goals / Setting high expectations"
$ws.Range("L1").Value = "Confidence in goals"
